# Update cached market-board / leve profit figures across sheets
# (scheduled runner refresh of Phantom_Profits data)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 898.9
$ws.Range("I12").Value = 899.1111
$ws.Range("J12").Value = 897
$ws.Range("K12").Value = 899.1111
$ws.Range("L12").Value = 897
$ws.Range("M12").Value = -729.1111
$ws.Range("N12").Value = -1237
$ws.Range("H32").Value = 3055.1428
$ws.Range("I32").Value = 1879.6
$ws.Range("K32").Value = 1879.6
$ws.Range("M32").Value = -1553.6
$ws.Range("H51").Value = 9811.75
$ws.Range("I51").Value = 9082.333000000001
$ws.Range("J51").Value = 12000
$ws.Range("K51").Value = 9082.333000000001
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = -8598.333000000001
$ws.Range("N51").Value = -12968
$ws.Range("H62").Value = 1149.5
$ws.Range("I62").Value = 1149.5
$ws.Range("K62").Value = 1149.5
$ws.Range("M62").Value = -525.5
$ws.Range("H65").Value = 1149.5
$ws.Range("I65").Value = 1149.5
$ws.Range("K65").Value = 5747.5
$ws.Range("M65").Value = -2627.5
$ws.Range("H75").Value = 57500
$ws.Range("J75").Value = 57500
$ws.Range("L75").Value = 57500
$ws.Range("N75").Value = -59372
$ws.Range("H78").Value = 57500
$ws.Range("J78").Value = 57500
$ws.Range("L78").Value = 172500
$ws.Range("N78").Value = -181860
$ws.Range("H100").Value = 2712.8462
$ws.Range("I100").Value = 3291.25
$ws.Range("J100").Value = 1787.4
$ws.Range("K100").Value = 3291.25
$ws.Range("L100").Value = 1787.4
$ws.Range("M100").Value = -2750.25
$ws.Range("N100").Value = -2869.4
$ws.Range("H137").Value = 3874.389
$ws.Range("J137").Value = 3882.2
$ws.Range("L137").Value = 11646.6
$ws.Range("N137").Value = -16746.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 15270
$ws.Range("J11").Value = 15270
$ws.Range("L11").Value = 15270
$ws.Range("N11").Value = -15558
$ws.Range("H131").Value = 75249.5
$ws.Range("J131").Value = 75249.5
$ws.Range("L131").Value = 75249.5
$ws.Range("N131").Value = -85329.5
$ws.Range("H132").Value = 2618.818
$ws.Range("I132").Value = 2618.818
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7856.454000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -5326.454000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 512
$ws.Range("I22").Value = 512
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 512
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -339
$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35630
$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37184
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H134").Value = 5736.5625
$ws.Range("I134").Value = 6264.3335
$ws.Range("K134").Value = 18793.0005
$ws.Range("M134").Value = -16258.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9999
$ws.Range("I86").Value = 9999
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 9999
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8876
$ws.Range("H88").Value = 6713
$ws.Range("J88").Value = 6713
$ws.Range("L88").Value = 6713
$ws.Range("N88").Value = -7525
$ws.Range("H89").Value = 9999
$ws.Range("I89").Value = 9999
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 49995
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -44379
$ws.Range("H91").Value = 6713
$ws.Range("J91").Value = 6713
$ws.Range("L91").Value = 6713
$ws.Range("N91").Value = -9521

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 122.833336
$ws.Range("J2").Value = 139
$ws.Range("L2").Value = 834
$ws.Range("N2").Value = -1060
$ws.Range("H5").Value = 1211.5
$ws.Range("J5").Value = 1489.5
$ws.Range("L5").Value = 4468.5
$ws.Range("N5").Value = -4692.5
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H129").Value = 2850.4
$ws.Range("I129").Value = 1999
$ws.Range("J129").Value = 3063.25
$ws.Range("K129").Value = 5997
$ws.Range("L129").Value = 9189.75
$ws.Range("M129").Value = -997
$ws.Range("N129").Value = -19189.75
$ws.Range("H135").Value = 1211.5
$ws.Range("J135").Value = 1489.5
$ws.Range("L135").Value = 13405.5
$ws.Range("N135").Value = -18475.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 269420
$ws.Range("I10").Value = 1000000
$ws.Range("J10").Value = 25893.334
$ws.Range("K10").Value = 1000000
$ws.Range("L10").Value = 25893.334
$ws.Range("M10").Value = -999831
$ws.Range("N10").Value = -26231.334
$ws.Range("H128").Value = 92563.39999999999
$ws.Range("J128").Value = 95456.75
$ws.Range("L128").Value = 95456.75
$ws.Range("N128").Value = -105416.75
$ws.Range("H132").Value = 6554.3335
$ws.Range("I132").Value = 4824.75
$ws.Range("K132").Value = 14474.25
$ws.Range("M132").Value = -11944.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 9321
$ws.Range("J76").Value = 9321
$ws.Range("L76").Value = 9321
$ws.Range("N76").Value = -9997
$ws.Range("H79").Value = 9321
$ws.Range("J79").Value = 9321
$ws.Range("L79").Value = 9321
$ws.Range("N79").Value = -11661
$ws.Range("H128").Value = 83980
$ws.Range("J128").Value = 83980
$ws.Range("L128").Value = 83980
$ws.Range("N128").Value = -93940
$ws.Range("H130").Value = 66663.336
$ws.Range("J130").Value = 66663.336
$ws.Range("L130").Value = 66663.336
$ws.Range("N130").Value = -76703.336
$ws.Range("H132").Value = 6959.8
$ws.Range("I132").Value = 6959.8
$ws.Range("K132").Value = 20879.4
$ws.Range("M132").Value = -18349.4
$ws.Range("H136").Value = 3249.125
$ws.Range("J136").Value = 4999.5
$ws.Range("L136").Value = 14998.5
$ws.Range("N136").Value = -20098.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 44000
$ws.Range("J112").Value = 44000
$ws.Range("L112").Value = 44000
$ws.Range("N112").Value = -46954
$ws.Range("H130").Value = 39332.668
$ws.Range("J130").Value = 39332.668
$ws.Range("L130").Value = 39332.668
$ws.Range("N130").Value = -49372.668
$ws.Range("H132").Value = 1758.2894
$ws.Range("I132").Value = 1797.8649
$ws.Range("K132").Value = 5393.5947
$ws.Range("M132").Value = -2863.5947
